$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.401.00'
$ws.Range("E2").Value = '  +2.89%  '
$ws.Range("D3").Value = '3.432.87'
$ws.Range("E3").Value = '  +2.11%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '406.73'
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.07'
$ws.Range("E6").Value = '  +4.67%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.596'
$ws.Range("E7").Value = '  -0.85%  '
$ws.Range("E9").Value = '  +3.65%  '
$ws.Range("E10").Value = '  +8.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.97'
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.86'
$ws.Range("E13").Value = '  +1.90%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.40'
$ws.Range("E14").Value = '  -0.54%  '
$ws.Range("D15").Value = '3.430.94'
$ws.Range("E15").Value = '  +2.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '11.69'
$ws.Range("E16").Value = '  +3.25%  '
$ws.Range("D17").Value = '62.258.46'
$ws.Range("E17").Value = '  +2.54%  '
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000148'
$ws.Range("E19").Value = '  +11.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.16'
$ws.Range("E20").Value = '  -2.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '84.18'
$ws.Range("E21").Value = '  +2.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '311.71'
$ws.Range("E22").Value = '  +2.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.75'
$ws.Range("E23").Value = '  -2.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.17'
$ws.Range("E24").Value = '  +1.29%  '
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '29.74'
$ws.Range("E26").Value = '  +1.16%  '
$ws.Range("E27").Value = '  -3.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.81'
$ws.Range("E28").Value = '  +4.97%  '
$ws.Range("E29").Value = '  +6.36%  '
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '43.92'
$ws.Range("E31").Value = '  +3.57%  '
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.34'
$ws.Range("E33").Value = '  -3.03%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("E35").Value = '  +0.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.62'
$ws.Range("E36").Value = '  -0.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("E38").Value = '  +1.58%  '
$ws.Range("E39").Value = '  -1.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.317'
$ws.Range("E40").Value = '  +12.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '142.29'
$ws.Range("E41").Value = '  +4.30%  '
$ws.Range("E42").Value = '  +0.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.98'
$ws.Range("E43").Value = '  -2.36%  '
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.80'
$ws.Range("E45").Value = '  -0.34%  '
$ws.Range("E46").Value = '  +0.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.26'
$ws.Range("E47").Value = '  -2.63%  '
$ws.Range("D48").Value = '2.104.88'
$ws.Range("E48").Value = '  -1.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.31'
$ws.Range("E49").Value = '  -1.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.94'
$ws.Range("E50").Value = '  +2.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.73'
$ws.Range("E51").Value = '  +20.73%  '